$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad) rows 2-189: update date serial from 45189 (2023-09-20)
# to 45190 (2023-09-21) for every data row.
for ($row = 2; $row -le 189; $row++) {
    $ws.Cells.Item($row, 3).Value = 45190
}
